$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- New row 0 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>153</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>SetInterval</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve">method, time in </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>ms</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>)</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Set interval calls the corresponding method. Here we aren’t the ones calling the method that’s why we are not putting any parenthesis here.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>This is called a higher order function.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>If we write a function inside a higher order function like,</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>setInterval</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>function(</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t>){</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>},2000);</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 1 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>155</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Arrays</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t xml:space="preserve">Same as java. However, </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>It supports any type of data inside.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">it supports dynamic addition. 2 </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>legnth</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> array. Want to add a new item?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>Array[</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t>2] = “mango”;</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Array is now of length 3.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>For more info go to MDN</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 2 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>156</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Array addition</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Array.push</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(“content”) &lt;= same as list</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Array.p</w:t>
            </w:r>
            <w:r>
              <w:t>op</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(“content”) &lt;= same as list</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Array.shift</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(); &lt;= pops the leftmost element and returns it.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Array.</w:t>
            </w:r>
            <w:r>
              <w:t>un</w:t>
            </w:r>
            <w:r>
              <w:t>shift</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(</w:t>
            </w:r>
            <w:r>
              <w:t>“value”</w:t>
            </w:r>
            <w:r>
              <w:t>); &lt;= p</w:t>
            </w:r>
            <w:r>
              <w:t>ushes</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> the </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">value in the </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">leftmost element </w:t>
            </w:r>
            <w:r>
              <w:t>(</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>array[</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t>0])</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 3 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>163</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>forEach</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(parameters)</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>If you only give one argument, then it's the element, every time.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>If you give two, then it's the element, then the index, in that order.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>If you give three, then it's element, index, and array, in that order.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>Ex:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Array.forEach</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>(a</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>){</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Console.log(a);</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>}</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 4 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>166</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Function name VS function name ()</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t xml:space="preserve">Just calling the function wouldn’t </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>execte</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> the code however when we say </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>function(</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve">)then the code is </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>ececuted</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>.</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 5 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>166</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Creating our own method on array</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>Array.prototype.functionName</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> = function(parameter){</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Function content</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>//this refers to the Array on which the method/function is called on.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>}</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()

# --- New row 6 ---
$newRow = $t.Rows.Add()
$cell = $newRow.Cells.Item(1)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>167</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(2)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>JS objects</w:t>
            </w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
$cell = $newRow.Cells.Item(3)
$p = $cell.Range.Paragraphs.Item(1)
$ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r>
              <w:t>Array isn’t the best option all the time. We can use JS objects which stores values in a key-value pair.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Var person </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>={</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
          </w:p>
          <w:p>
            <w:r>
              <w:t>name: “”,</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>id: 21,</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>city: “NYC”</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>}</w:t>
</w:r>
          </w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$ip.InsertXML($xml)
$cell.Range.Paragraphs.Item(1).Range.Delete()
